$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update emphasized cell values (append footnote markers) ---
# Order matters for shared-string table ordering: base-dagger, then
# Stationary-dagger, then Gluconate-asterisks, then Glucose-asterisks,
# then high-asterisks.
$dagger = [char]0x2020

$ws.Range("B5").Value  = "base" + $dagger
$ws.Range("E7").Value  = "Stationary" + $dagger
$ws.Range("E16").Value = "Stationary" + $dagger

$ws.Range("D12").Value = "Gluconate***"
$ws.Range("D13").Value = "Gluconate***"

$ws.Range("D4").Value  = "Glucose***"

$ws.Range("C3").Value = "high***"
$ws.Range("C4").Value = "high***"
$ws.Range("C5").Value = "high***"

# --- Reset the special emphasis colors (red / accent blue) back to the
# normal automatic/black text color now that the emphasis is conveyed by
# the appended markers instead of font color. ---
$markedCells = @("C3","C4","D4","B5","C5","E7","D12","D13","E16")
foreach ($addr in $markedCells) {
    $ws.Range($addr).Font.ThemeColor = 1
}

# --- Selection moved from A1:E16 to A11:E16, active cell E16 ---
$ws.Range("A11:E16").Select
$excel.ActiveWindow.RangeSelection.Item(1).Activate() | Out-Null
